# Build v4.22.00.158 - update Employee roster with the new
# "RTGO Operator 2023-12-07" role group and its associated employees.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$ws.Range("A3").Value = "ROLE GROUP : RTGO Operator 2023-12-06T09:13:59.345"
$ws.Range("A4").Value = "92457737 - Lilliana Williamson`nROLE : RTGO100 1701844270281"
$ws.Range("A5").Value = "90317880 - Lewis Mosciski`nROLE : RTGO100 1701844270281"
$ws.Range("A6").Value = "90833312 - Angelo Mueller`nROLE : RTGO100 1701844270281"
$ws.Range("A7").Value = "ROLE GROUP : RTGO Operator 2023-12-07T19:27:58.156908600"
$ws.Range("A8").Value = "92970163 - Glenna Lynch`nROLE : RTGO100 1701853905917"

# Restore default row heights after writing the multi-line values so the
# rows don't end up flagged with an explicit custom height.
$ws.Rows("3:8").AutoFit()
